# Auto-generated Excel COM-interop script to apply numeric updates
# to the Leve profit calculation columns (H-N) across all 8 sheets,
# per the scheduled runner's refreshed market price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 182.5
$ws.Range("I2").Value = 160
$ws.Range("K2").Value = 160
$ws.Range("M2").Value = -47

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 125269
$ws.Range("I9").Value = 200150.4
$ws.Range("K9").Value = 200150.4
$ws.Range("M9").Value = -199981.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 404.57144
$ws.Range("J41").Value = 565.6667
$ws.Range("L41").Value = 565.6667
$ws.Range("N41").Value = -1445.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 900
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 659.1852
$ws.Range("I92").Value = 716.75
$ws.Range("J92").Value = 494.7143
$ws.Range("K92").Value = 716.75
$ws.Range("L92").Value = 494.7143
$ws.Range("M92").Value = 531.25
$ws.Range("N92").Value = -2990.7143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1431.75
$ws.Range("J97").Value = 1431.75
$ws.Range("L97").Value = 4295.25
$ws.Range("N97").Value = -5287.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2707.1765
$ws.Range("I100").Value = 2813.875
$ws.Range("K100").Value = 2813.875
$ws.Range("M100").Value = -2272.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1681
$ws.Range("I103").Value = 949.5
$ws.Range("J103").Value = 1863.875
$ws.Range("K103").Value = 2848.5
$ws.Range("L103").Value = 5591.625
$ws.Range("M103").Value = -2262.5
$ws.Range("N103").Value = -6763.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1366.2667
$ws.Range("I107").Value = 1508.4166
$ws.Range("K107").Value = 1508.4166
$ws.Range("M107").Value = 411.5834

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5550
$ws.Range("I116").Value = 5550
$ws.Range("K116").Value = 5550
$ws.Range("M116").Value = -2108

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 933.2
$ws.Range("I132").Value = 933.2
$ws.Range("K132").Value = 2799.6
$ws.Range("M132").Value = -269.6000000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4111.0557
$ws.Range("J138").Value = 4872.3335
$ws.Range("L138").Value = 14617.0005
$ws.Range("N138").Value = -24897.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1543.5217
$ws.Range("I2").Value = 1500.1578
$ws.Range("K2").Value = 1500.1578
$ws.Range("M2").Value = -1387.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22148
$ws.Range("I74").Value = 22922.857
$ws.Range("K74").Value = 22922.857
$ws.Range("M74").Value = -22048.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 22148
$ws.Range("I77").Value = 22922.857
$ws.Range("K77").Value = 114614.285
$ws.Range("M77").Value = -110246.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1845.5555
$ws.Range("I97").Value = 479.08334
$ws.Range("J97").Value = 4578.5
$ws.Range("K97").Value = 479.08334
$ws.Range("L97").Value = 4578.5
$ws.Range("M97").Value = 16.91665999999998
$ws.Range("N97").Value = -5570.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1122.3334
$ws.Range("I102").Value = 1122.3334
$ws.Range("K102").Value = 1122.3334
$ws.Range("M102").Value = 499.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1543.5217
$ws.Range("I116").Value = 1500.1578
$ws.Range("K116").Value = 1500.1578
$ws.Range("M116").Value = 793.8422

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3582.25
$ws.Range("I122").Value = 2998
$ws.Range("K122").Value = 8994
$ws.Range("M122").Value = -6544

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1543.5217
$ws.Range("I3").Value = 1500.1578
$ws.Range("K3").Value = 1500.1578
$ws.Range("M3").Value = -1386.1578

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 150
$ws.Range("J22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("N22").Value = -596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2588.5
$ws.Range("I94").Value = 3575.5
$ws.Range("J94").Value = 943.5
$ws.Range("K94").Value = 3575.5
$ws.Range("L94").Value = 943.5
$ws.Range("M94").Value = -3124.5
$ws.Range("N94").Value = -1845.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 208.6
$ws.Range("I7").Value = 191.66667
$ws.Range("J7").Value = 234
$ws.Range("K7").Value = 191.66667
$ws.Range("L7").Value = 234
$ws.Range("M7").Value = -78.66667000000001
$ws.Range("N7").Value = -460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2968.8
$ws.Range("I16").Value = 1450
$ws.Range("J16").Value = 3981.3333
$ws.Range("K16").Value = 1450
$ws.Range("L16").Value = 3981.3333
$ws.Range("M16").Value = -1163
$ws.Range("N16").Value = -4555.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1912.1818
$ws.Range("I31").Value = 1455.2858
$ws.Range("J31").Value = 2711.75
$ws.Range("K31").Value = 1455.2858
$ws.Range("L31").Value = 2711.75
$ws.Range("M31").Value = -1160.2858
$ws.Range("N31").Value = -3301.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1912.1818
$ws.Range("I34").Value = 1455.2858
$ws.Range("J34").Value = 2711.75
$ws.Range("K34").Value = 1455.2858
$ws.Range("L34").Value = 2711.75
$ws.Range("M34").Value = -1253.2858
$ws.Range("N34").Value = -3115.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 30000
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 43042
$ws.Range("J50").Value = 43042
$ws.Range("L50").Value = 43042
$ws.Range("N50").Value = -44292

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 14165
$ws.Range("J59").Value = 14165
$ws.Range("L59").Value = 14165
$ws.Range("N59").Value = -16455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69995
$ws.Range("J68").Value = 69995
$ws.Range("L68").Value = 69995
$ws.Range("N68").Value = -71493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 69995
$ws.Range("J71").Value = 69995
$ws.Range("L71").Value = 209985
$ws.Range("N71").Value = -217473

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2968.8
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 3981.3333
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 3981.3333
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -8321.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5141.125
$ws.Range("I122").Value = 4559.8335
$ws.Range("J122").Value = 6885
$ws.Range("K122").Value = 13679.5005
$ws.Range("L122").Value = 20655
$ws.Range("M122").Value = -11229.5005
$ws.Range("N122").Value = -25555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2388
$ws.Range("J46").Value = 2388
$ws.Range("L46").Value = 7164
$ws.Range("N46").Value = -7346

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2995.2964
$ws.Range("J122").Value = 3421.1738
$ws.Range("L122").Value = 30790.5642
$ws.Range("N122").Value = -35690.5642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2610.3333
$ws.Range("J137").Value = 2998.5
$ws.Range("L137").Value = 8995.5
$ws.Range("N137").Value = -19195.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 707.6667
$ws.Range("I102").Value = 649.2
$ws.Range("K102").Value = 649.2
$ws.Range("M102").Value = 972.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2633.0833
$ws.Range("I122").Value = 2259.7
$ws.Range("K122").Value = 6779.099999999999
$ws.Range("M122").Value = -4329.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2466
$ws.Range("I126").Value = 2466
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7398
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4928
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2099.1667
$ws.Range("I46").Value = 973.75
$ws.Range("J46").Value = 4350
$ws.Range("K46").Value = 973.75
$ws.Range("L46").Value = 4350
$ws.Range("M46").Value = -785.75
$ws.Range("N46").Value = -4726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1823.2174
$ws.Range("I93").Value = 1525.8
$ws.Range("J93").Value = 2380.875
$ws.Range("K93").Value = 1525.8
$ws.Range("L93").Value = 2380.875
$ws.Range("M93").Value = -277.8
$ws.Range("N93").Value = -4876.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5749.25
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4999
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5393
$ws.Range("I136").Value = 4524.2
$ws.Range("J136").Value = 7999.4
$ws.Range("K136").Value = 13572.6
$ws.Range("L136").Value = 23998.2
$ws.Range("M136").Value = -11022.6
$ws.Range("N136").Value = -29098.2
